# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock)
# worksheet, matching the new normalized-export schema.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票

# Current layout (before):
#   A        B     C      D         E           F         G      H                    I       J               K
#   (idx)    name  owner  quantity  face_value  currency  total  property_category   date    legislator_name legislator_id
#
# Target layout (after):
#   A        B     C      D         E           F         G      H                    I         J       K               L              M            N
#   (idx)    name  owner  quantity  face_value  currency  total  property_category   category  date    legislator_name legislator_id  source_file  index

$lastRow = 39

# Insert a new blank column before the existing "date" column (column I),
# shifting date/legislator_name/legislator_id one column to the right.
$ws.Columns.Item(9).Insert()

# New column I is "category" -- every stock row belongs to the "normal" export.
$ws.Range("I1").Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# Append the two new trailing columns: source_file + index.
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmp700a1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
